$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the paired header columns from "<name>_old" / "<name>_new" to
# "<name>_FV2210" / "<name>_FV2304" respectively (columns A-J pair with L-U).
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $leftCol = $i + 1
    $rightCol = $i + 12
    $ws.Cells.Item(1, $leftCol).Value = "$($baseNames[$i])_FV2210"
    $ws.Cells.Item(1, $rightCol).Value = "$($baseNames[$i])_FV2304"
}

# Turn the used range into a real Excel table (adds the autofilter too).
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U66"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
